$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.361.09"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.687.32"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").Value = "'218.18"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'0.5455"
$ws.Range("E6").Value = "  +4.05%  "
$ws.Range("D8").Value = "'0.2721"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").Value = "'0.06460"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'21.99"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "'0.07684"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").Value = "1.683.85"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "'4.532"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "'0.5805"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "'0.000008384"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "'65.10"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "26.413.98"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "'10.97"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "'190.55"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'6.238"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D24").Value = "'149.72"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "'0.1310"
$ws.Range("E25").Value = "  +5.17%  "
$ws.Range("D26").Value = "'7.866"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").Value = "'15.70"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "'0.06306"
$ws.Range("E28").Value = "  -9.15%  "
$ws.Range("D29").Value = "'1.407"
$ws.Range("E29").Value = "  +4.84%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'3.586"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "'1.681"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("D35").Value = "'0.6159"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").Value = "'2.412"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'2.720"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").Value = "'6.260"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "1.112.80"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("D40").Value = "'0.01625"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "'0.8831"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'101.37"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "1.837.77"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.011"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.177"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "'0.05271"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'0.4307"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'6.036"
$ws.Range("E51").Value = "  +0.44%  "
